# Implemented getting kafka relations.
# The "classFields" sheet lists, per class, the fields discovered for that
# class (Field Name / Field Modifier / Field Type). Re-implementing field
# discovery to also pick up the kafka-related fields changed the order in
# which fields are enumerated for several classes. This reorders the rows
# for the affected classes to match the new enumeration order, leaving the
# header row and already-correct rows untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("classFields")

function Set-FieldRow($row, $name, $modifier, $type) {
    $ws.Cells.Item($row, 2).Value = $name
    $ws.Cells.Item($row, 3).Value = $modifier
    $ws.Cells.Item($row, 4).Value = $type
}

# pl.piomin.payment.service.OrderManageService (rows 2-5):
#   template, SOURCE, LOG, repository  ->  LOG, repository, template, SOURCE
Set-FieldRow 2 "LOG" "private" "org.slf4j.Logger"
Set-FieldRow 3 "repository" "private" "pl.piomin.payment.repository.CustomerRepository"
Set-FieldRow 4 "template" "private" "org.springframework.kafka.core.KafkaTemplate"
Set-FieldRow 5 "SOURCE" "private" "java.lang.String"

# pl.piomin.payment.PaymentComponentTests (rows 6-11):
#   LOG, factory, customer, kafka, template, repository
#   -> template, customer, repository, kafka, LOG, factory
Set-FieldRow 6 "template" "private" "org.springframework.kafka.core.KafkaTemplate"
Set-FieldRow 7 "customer" "" "pl.piomin.payment.domain.Customer"
Set-FieldRow 8 "repository" "" "pl.piomin.payment.repository.CustomerRepository"
Set-FieldRow 9 "kafka" "private" "org.springframework.kafka.test.EmbeddedKafkaBroker"
Set-FieldRow 10 "LOG" "private" "org.slf4j.Logger"
Set-FieldRow 11 "factory" "private" "org.springframework.kafka.core.ConsumerFactory"

# pl.piomin.payment.PaymentApp (rows 12-14): order unchanged.

# pl.piomin.payment.domain.Customer (rows 15-18):
#   id, name, amountReserved, amountAvailable
#   -> name, amountReserved, id, amountAvailable
Set-FieldRow 15 "name" "private" "java.lang.String"
Set-FieldRow 16 "amountReserved" "private" "int"
Set-FieldRow 17 "id" "private" "java.lang.Long"
# row 18 (amountAvailable) unchanged.
